$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '58.002.08'
$ws.Range("E2").Value = '  +0.52%  '

# Row 3
$ws.Range("D3").Value = '2.478.96'
$ws.Range("E3").Value = '  +1.29%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.31%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '516.99'
$ws.Range("E5").Value = '  +0.36%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.50'
$ws.Range("E6").Value = '  -0.56%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.31%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.552'
$ws.Range("E8").Value = '  -0.37%  '

# Row 9
$ws.Range("D9").Value = '2.492.17'
$ws.Range("E9").Value = '  +1.78%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0963'
$ws.Range("E10").Value = '  -1.26%  '

# Row 11
$ws.Range("E11").Value = '  -0.29%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.10'
$ws.Range("E12").Value = '  -2.54%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.327'
$ws.Range("E13").Value = '  -3.07%  '

# Row 14
$ws.Range("D14").Value = '2.917.79'
$ws.Range("E14").Value = '  +1.17%  '

# Row 15
$ws.Range("D15").Value = '57.938.28'
$ws.Range("E15").Value = '  +0.45%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.87'
$ws.Range("E16").Value = '  -0.87%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000133'
$ws.Range("E17").Value = '  -0.87%  '

# Row 18
$ws.Range("D18").Value = '2.491.35'
$ws.Range("E18").Value = '  +0.94%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.61'
$ws.Range("E19").Value = '  +0.04%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '318.94'
$ws.Range("E20").Value = '  +0.30%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.13'
$ws.Range("E21").Value = '  -0.29%  '

# Row 22
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.995'
$ws.Range("E22").Value = '  -0.32%  '

# Row 23
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.94'
$ws.Range("E23").Value = '  +4.44%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.62'
$ws.Range("E24").Value = '  -0.57%  '

# Row 25
$ws.Range("B25").Value = 'Binance-PegBSC-USD'
$ws.Range("C25").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.995'
$ws.Range("E25").Value = '  -0.74%  '

# Row 26
$ws.Range("B26").Value = 'Polygon'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.397'
$ws.Range("E26").Value = '  -1.86%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.159'
$ws.Range("E27").Value = '  +0.21%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.25'
$ws.Range("E28").Value = '  -0.21%  '

# Row 29
$ws.Range("D29").Value = '0.0₃0741'
$ws.Range("E29").Value = '  +1.24%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '167.33'
$ws.Range("E30").Value = '  +1.06%  '

# Row 31
$ws.Range("E31").Value = '  +0.73%  '

# Row 32
$ws.Range("E32").Value = '  +2.17%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.13'
$ws.Range("E33").Value = '  -0.86%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.998'
$ws.Range("E34").Value = '  -0.10%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.995'
$ws.Range("E35").Value = '  -0.33%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.90'
$ws.Range("E36").Value = '  -0.16%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.24'
$ws.Range("E37").Value = '  -3.86%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.89'

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.62'
$ws.Range("E39").Value = '  +1.74%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.44'
$ws.Range("E40").Value = '  -1.29%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.763'
$ws.Range("E41").Value = '  -2.40%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '272.39'
$ws.Range("E42").Value = '  +0.83%  '

# Row 43
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.02'
$ws.Range("E43").Value = '  +1.42%  '

# Row 44
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.38'
$ws.Range("E44").Value = '  -0.53%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.590'
$ws.Range("E45").Value = '  +0.75%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0913'
$ws.Range("E46").Value = '  +1.04%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '120.54'
$ws.Range("E47").Value = '  -2.23%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0495'
$ws.Range("E48").Value = '  +2.39%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '17.56'
$ws.Range("E49").Value = '  +0.32%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0210'
$ws.Range("E50").Value = '  +0.95%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.61'
$ws.Range("E51").Value = '  +0.02%  '
